# Progress.xlsx update
# - Fill in previously-empty rows 8-15 with parsed-range statistics
# - Fix C2/D2 values
# - Re-style columns A, B (italic font, top-aligned) and E (wrap text)
# - Resize rows 8 and 12 to fit wrapped content
# - Update the saved selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix existing values in row 2 (C2, D2)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 996
$ws.Range("D2").Value = 23

# ---------------------------------------------------------------------------
# 2) Fill in the previously-blank data rows (8-15)
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = 932
$ws.Range("D8").Value = 14
$ws.Range("E8").Value = "6059, 6122, 6123, 6124, 6125, 6126, 6127, 6128, 6129, 6130, 6131, 6132, 6133, 6134, 6135, 6136, 6137, 6138, 6139, 6140, 6141, 6142, 6143, 6144, 6145, 6146, 6147, 6148, 6149, 6150, 6151, 6152, 6153, 6154, 6155, 6156, 6157, 6158, 6159, 6160, 6161, 6162, 6163, 6164, 6165, 6440, 6517, 6518, 6519, 6520, 6642, 6643, 6644, 6645, 6646, 6647, 6648, 6649, 6650, 6651, 6652, 6653, 6654, 6655, 6656, 6657, 6658, 6728"

$ws.Range("C9").Value = 999
$ws.Range("D9").Value = 14
$ws.Range("E9").Value = 7236

$ws.Range("C10").Value = 996
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "8073, 8086, 8706, 8797, 8879"

$ws.Range("C11").Value = 996
$ws.Range("D11").Value = 7
$ws.Range("E11").Value = "-"

$ws.Range("C12").Value = 982
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "10062, 10070, 10127, 10432, 10436, 10441, 10463, 10470, 10501, 10577, 10584, 10589, 10595, 10609, 10671, 10761, 10762, 10888, 10959"

$ws.Range("C13").Value = 992
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "11099, 11306, 11307, 11358, 11384, 11427, 11428, 11709, 11990"

$ws.Range("C14").Value = 995
$ws.Range("D14").Value = 8
$ws.Range("E14").Value = "12014, 12279, 12295, 12417, 12579"

$ws.Range("B15").Value = 13520
$ws.Range("C15").Value = 513
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = "13061, 13123, 13125, 13209, 13217, 13465, 13496"

# ---------------------------------------------------------------------------
# 3) Re-style column A (id "from") and column B (id "to"): italic font,
#    vertical alignment changes from center to top.
# ---------------------------------------------------------------------------
$colA = $ws.Range("A2:A15")
$colA.VerticalAlignment = -4160
$colA.HorizontalAlignment = -4152
$colA.Font.Italic = $true

$colB = $ws.Range("B2:B15")
$colB.VerticalAlignment = -4160
$colB.HorizontalAlignment = -4131
$colB.Font.Italic = $true

# ---------------------------------------------------------------------------
# 4) Column E (and D13, which now carries a text dash instead of a number)
#    wraps long comma-separated id lists.
# ---------------------------------------------------------------------------
$colE = $ws.Range("E2:E15")
$colE.WrapText = $true

$ws.Range("D13").WrapText = $true

# ---------------------------------------------------------------------------
# 5) Resize rows that now contain wrapped multi-line content.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 90
$ws.Rows.Item(12).RowHeight = 30

# ---------------------------------------------------------------------------
# 6) Restore the last selected cell.
# ---------------------------------------------------------------------------
$ws.Range("D24").Select()
